$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Cells.Item(1, 1).Value = "Kim"
$ws.Cells.Item(1, 2).Value = "Point"
$ws.Cells.Item(1, 3).Value = "Emil"
$ws.Cells.Item(1, 4).Value = "Point"
$ws.Cells.Item(1, 5).Value = "Mads"
$ws.Cells.Item(1, 6).Value = "Point"
$ws.Cells.Item(1, 7).Value = "Soren"
$ws.Cells.Item(1, 8).Value = "Point"

# Matchup rows
$ws.Cells.Item(2, 1).Value = "Chelsea"
$ws.Cells.Item(2, 2).Value = 0
$ws.Cells.Item(2, 3).Value = "Dortmund"
$ws.Cells.Item(2, 4).Value = 0
$ws.Cells.Item(2, 5).Value = "Ac Milan"
$ws.Cells.Item(2, 6).Value = 0
$ws.Cells.Item(2, 7).Value = "Fc midtjylland"
$ws.Cells.Item(2, 8).Value = 0

$ws.Cells.Item(3, 1).Value = "Barcelona"
$ws.Cells.Item(3, 2).Value = 0
$ws.Cells.Item(3, 3).Value = "Atalanta"
$ws.Cells.Item(3, 4).Value = 0
$ws.Cells.Item(3, 5).Value = "FC Kbenhavn"
$ws.Cells.Item(3, 6).Value = 0
$ws.Cells.Item(3, 7).Value = "sevilla"
$ws.Cells.Item(3, 8).Value = 0

$ws.Cells.Item(4, 1).Value = "RB Leipzig"
$ws.Cells.Item(4, 2).Value = 0
$ws.Cells.Item(4, 3).Value = "Real sociedad"
$ws.Cells.Item(4, 4).Value = 0
$ws.Cells.Item(4, 5).Value = "Frankfurt"
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = "juventus"
$ws.Cells.Item(4, 8).Value = 0

$ws.Cells.Item(5, 1).Value = "AGF"
$ws.Cells.Item(5, 2).Value = 0
$ws.Cells.Item(5, 3).Value = "Tottenham"
$ws.Cells.Item(5, 4).Value = 0
$ws.Cells.Item(5, 5).Value = "Arsenal"
$ws.Cells.Item(5, 6).Value = 0
$ws.Cells.Item(5, 7).Value = "Manchester Utd"
$ws.Cells.Item(5, 8).Value = 0

$ws.Cells.Item(6, 1).Value = "Bologna"
$ws.Cells.Item(6, 2).Value = 0
$ws.Cells.Item(6, 3).Value = "Brndby IF"
$ws.Cells.Item(6, 4).Value = 0
$ws.Cells.Item(6, 5).Value = "Valencia"
$ws.Cells.Item(6, 6).Value = 0
$ws.Cells.Item(6, 7).Value = "Leverkusen"
$ws.Cells.Item(6, 8).Value = 0

$ws.Cells.Item(7, 1).Value = "OB"
$ws.Cells.Item(7, 2).Value = 0
$ws.Cells.Item(7, 3).Value = "Leicester"
$ws.Cells.Item(7, 4).Value = 0
$ws.Cells.Item(7, 5).Value = "Hoffenheim"
$ws.Cells.Item(7, 6).Value = 0
$ws.Cells.Item(7, 7).Value = "Torino"
$ws.Cells.Item(7, 8).Value = 0

# Totals row
$ws.Cells.Item(8, 1).Value = "Total:"
$ws.Cells.Item(8, 2).Formula = "=SUM(B2:B7)"
$ws.Cells.Item(8, 3).Value = "Total:"
$ws.Cells.Item(8, 4).Formula = "=SUM(D2:D7)"
$ws.Cells.Item(8, 5).Value = "Total:"
$ws.Cells.Item(8, 6).Formula = "=SUM(F2:F7)"
$ws.Cells.Item(8, 7).Value = "Total:"
$ws.Cells.Item(8, 8).Formula = "=SUM(H2:H7)"
